# Refactor: rename the "activity" concept to "step" throughout the
# workbook (sheet names + id/name columns), and record two new
# "Traject" step-property rows (categories A / Q) on the property sheet.
#
# NOTE: operations below are ordered deliberately so that brand-new
# strings are first written to the workbook in the same sequence the
# original author produced them in (shared-string table entries are
# assigned in first-write order), and so the sheet that ends up
# active/selected matches the original file ("step", née "activity").

$wb = $excel.ActiveWorkbook

$wsStep     = $wb.Worksheets.Item("activity")
$wsFlow     = $wb.Worksheets.Item("process_flow")
$wsStepProp = $wb.Worksheets.Item("activity_property")

# --- "activity_property" sheet: two new rows for a "Traject" property ---
$wsStepProp.Range("A5").Value = 8
$wsStepProp.Range("B5").Value = "Traject"
$wsStepProp.Range("C5").Value = "A"
$wsStepProp.Range("D5").Value = 0.3

$wsStepProp.Range("A6").Value = 8
$wsStepProp.Range("B6").Value = "Traject"
$wsStepProp.Range("C6").Value = "Q"
$wsStepProp.Range("D6").Value = 0.7

$wsStepProp.PageSetup.PaperSize = 9      # A4
$wsStepProp.PageSetup.Orientation = 1    # xlPortrait

$wsStepProp.Name = "step_property"

# --- "activity" sheet -> "step": rename columns, use sentinel strings ---
$wsStep.Name = "step"

$wsStep.Range("A1").Value = "step_id"
$wsStep.Range("B1").Value = "step_name"

$wsStep.Range("A2").Value = "START"
$wsStep.Range("A3").Value = "END"

# --- "process_flow" sheet: same id rename + START/END sentinels ---
$wsFlow.Range("A1").Value = "step_id"
$wsFlow.Range("B1").Value = "next_step_id"

$wsFlow.Range("A2").Value = "START"
$wsFlow.Range("A3").Value = "START"

$wsFlow.Range("B11").Value = "END"
$wsFlow.Range("B14").Value = "END"

$wsFlow.Range("B14").Select() | Out-Null

# --- finish renaming "activity_property" -> "step_property" ---
$wsStepProp.Range("A1").Value = "step_id"

# leave the "step" sheet active, mirroring the saved file
$wsStep.Range("B2").Select() | Out-Null
